# edit.ps1 - applies the "escape_tpl.docx" template update:
#   - merge the 3 runs of the "Here a listing..." paragraph into one run
#   - drop the trailing space in the "{{r page_break }} " run
#   - insert a new "{{ new_listing }}" paragraph (orange) before "END"
#   - Normal style: set Font.Kerning to 0 (adds <w:kern w:val="0"/>)
#   - Titre style:  (re)assert KeepWithNext (keepNext) on the paragraph format

$d = $word.ActiveDocument

# --- 1. Merge the 3 runs that make up the "Here a listing ..." paragraph
#        into a single run by doing a same-text Find & Replace across them.
$mergedText = "Here a listing that escapes and manages newline and page break AND keep the current character styling :"
$d.Content.Find.Execute($mergedText, $true, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2) | Out-Null

# --- 2. Remove the trailing space after the page_break tag.
$d.Content.Find.Execute("{{r page_break }} ", $true, $false, $false, $false, $false, $true, 1, $false, "{{r page_break }}", 2) | Out-Null

# --- 3. Insert a new paragraph "{{ new_listing }}" right before the "END"
#        paragraph, styled with an orange font color (RGB C06616).
$endPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($endPara.Range.Start, $endPara.Range.Start)
$insertionPoint.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newRange = $newPara.Range
$newRange.Text = "{{ new_listing }}"

# Word's Font.Color is a BGR-packed long (wdColor), not plain RGB, so
# swap the byte order of 0xC06616 before assigning it.
$rr = 0xC0
$gg = 0x66
$bb = 0x16
$bgrColor = ($bb * 0x10000) + ($gg * 0x100) + $rr
$newRange.Font.Color = $bgrColor

# --- 4. Normal style: enable kerning at size 0 -> <w:kern w:val="0"/>
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.Font.Kerning = 0

# --- 5. Titre style: make sure the paragraph is kept with the next one.
$titreStyle = $d.Styles.Item("Titre")
$titreStyle.ParagraphFormat.KeepWithNext = $true
